$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Ngày công: 17 -> 19
$ws.Range("B2").Value = 19

# Phụ cấp: 595000 -> 665000
$ws.Range("B3").Value = 665000

# Lương cơ bản tại SÓC TRĂNG: 1821428.571428571 -> 2035714.285714286
$ws.Range("B20").Value = 2035714.285714286

# Tổng lương tại CẦN THƠ: 845000 -> 915000
$ws.Range("B28").Value = 915000

# Tổng lương tại SÓC TRĂNG: 1890228.571428571 -> 2104514.285714285
$ws.Range("B30").Value = 2104514.285714285

# Tổng lương: 2735228.571428571 -> 3019514.285714285
$ws.Range("B31").Value = 3019514.285714285
